# Update gh-pages output: new "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini" entry
# (2024-10-03) inserted, plus refreshed "want to go" counters (column F) on a
# handful of existing rows across the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

function Set-IndexCellStyle {
    param($ws, [int]$row)
    # Column A carries the bold/centered/bordered "index" style (style id 1
    # in the original workbook). Rows.Insert() on a fresh row doesn't always
    # pick that up, so stamp it explicitly to match the sibling index cells.
    $cell = $ws.Cells.Item($row, 1)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

function Set-DateTextCell {
    param($ws, [string]$addr, [string]$text)
    # Column B holds plain "yyyy-mm-dd" strings, not real dates. Force the
    # cell to Text first so Excel doesn't auto-convert the literal into a
    # date serial number.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Refresh "want to go" counters on existing rows.
$ws1.Range("F2").Value = 3326
$ws1.Range("F4").Value = 58
$ws1.Range("F5").Value = 1360

# Insert a new row 6 (pushes the "万圣漫控嘉年华10" row from 6 -> 7).
$ws1.Rows.Item(6).Insert()

$ws1.Range("A6").Value = 5
Set-DateTextCell $ws1 "B6" "2024-10-03"
$ws1.Range("C6").Value = "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini"
$ws1.Range("D6").Value = "南宁国际会展中心  南宁国际会展中心"
$ws1.Range("E6").Value = "2024.10.03 09:30-10.04 17:30"
$ws1.Range("F6").Value = 2
$ws1.Range("G6").Value = 55
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=91043"
$ws1.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202408/jEAI96Ev1724123680899.jpeg"
Set-IndexCellStyle $ws1 6

# The shifted-down row (now row 7) needs its index re-numbered and its
# "want to go" counter refreshed.
$ws1.Range("A7").Value = 6
$ws1.Range("F7").Value = 317

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types combined)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 3326
$ws4.Range("F4").Value = 58
$ws4.Range("F5").Value = 1360

# Insert a new row 6 (pushes "莫西干人..." row 6 -> 7, and "万圣漫控嘉年华10"
# row 7 -> 8).
$ws4.Rows.Item(6).Insert()

$ws4.Range("A6").Value = 5
Set-DateTextCell $ws4 "B6" "2024-10-03"
$ws4.Range("C6").Value = "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini"
$ws4.Range("D6").Value = "南宁国际会展中心  南宁国际会展中心"
$ws4.Range("E6").Value = "2024.10.03 09:30-10.04 17:30"
$ws4.Range("F6").Value = 2
$ws4.Range("G6").Value = 55
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=91043"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202408/jEAI96Ev1724123680899.jpeg"
Set-IndexCellStyle $ws4 6

# Shifted-down rows need re-numbered indices; the last one ("万圣漫控嘉年华10")
# also gets its "want to go" counter refreshed.
$ws4.Range("A7").Value = 6
$ws4.Range("A8").Value = 7
$ws4.Range("F8").Value = 317
